# Update the "Improved Water" contributing-factor label to
# "Improved Water Source" on every sheet that shows it (K2 on each of
# the three sheets shares the same text via the shared-strings table).
$wb = $excel.ActiveWorkbook

$instructions = $wb.Worksheets.Item("instructions")
$data = $wb.Worksheets.Item("Data")
$cat = $wb.Worksheets.Item("Cat")

$instructions.Range("K2").Value = "Improved Water Source"
$data.Range("K2").Value = "Improved Water Source"
$cat.Range("K2").Value = "Improved Water Source"

# Restore each sheet's active selection to match the saved view state.
$instructions.Activate()
$instructions.Range("H9").Select()

$data.Activate()
$data.Range("K3").Select()

$cat.Activate()
$cat.Range("C21").Select()

# "instructions" is the tab that was active/selected when the file was saved.
$instructions.Activate()
